# Appends four new list-item paragraphs (ilvl=1, numId=1) at the end of the
# document, after the "...duplicatas de objetos..." paragraph, mirroring the
# paragraph/run structure newly created by Word for the "Tuplas" notes.
#
# Strategy:
#  - Paragraph.Range.InsertParagraphAfter() on the last paragraph creates a
#    new paragraph that already inherits the same pStyle/numPr/spacing/jc and
#    the same run-level rPr (font/size) as the source paragraph, matching the
#    diff's new <w:p>/<w:pPr>/<w:rPr> blocks exactly.
#  - Setting the new paragraph's Range.Text in one shot writes a single run
#    with fully-inherited formatting (rFonts Times New Roman, sz/szCs 24).
#  - Where the diff shows that single logical sentence split across multiple
#    <w:r> elements (with identical rPr), we locate the split point with
#    Find (scoped to the new paragraph only) and toggle Bold on/off across
#    the sub-range; the no-op formatting change forces the run to be
#    serialized separately without altering the visible formatting.

$d = $word.ActiveDocument

function New-ListParagraph([string]$text) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Last
    $newP.Range.Text = $text
    return $newP
}

function Split-Run([object]$paragraph, [string]$anchorText) {
    # Splits the paragraph's run so that everything *before* the first
    # occurrence of $anchorText (scoped to this paragraph) becomes its own
    # run, separate from $anchorText onward.
    $scope = $d.Range($paragraph.Range.Start, $paragraph.Range.End)
    $scope.Find.ClearFormatting()
    $scope.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $before = $d.Range($paragraph.Range.Start, $scope.Start)
    $before.Bold = 1
    $before.Bold = 0
}

# --- Paragraph: " Uma tupla é como uma lista, porém imutável, ou seja, não
#     podemos usar append() ou remove()." -> split off the trailing "." ----
$p1 = New-ListParagraph(" Uma tupla é como uma lista, porém imutável, ou seja, não podemos usar append() ou remove().")
Split-Run $p1 "."

# --- Paragraph: " Podemos colocar valores de vários tipos: str, int, float,
#     dentre outros em uma tupla." -> single run, no split ----------------
$p2 = New-ListParagraph(" Podemos colocar valores de vários tipos: str, int, float, dentre outros em uma tupla.")

# --- Paragraph: " Diferentemente da lista, que é demarcada com [], as
#     tuplas são com ()." -> split off the leading " " -------------------
$p3 = New-ListParagraph(" Diferentemente da lista, que é demarcada com [], as tuplas são com ().")
Split-Run $p3 "Diferentemente"

# --- Paragraph: " A posição dos elementos também é significante,
#     diferentemente das listas." -> split into three runs ----------------
$p4 = New-ListParagraph(" A posição dos elementos também é significante, diferentemente das listas.")
Split-Run $p4 " listas."
Split-Run $p4 "listas."

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
